# GUI Addition and scraper cleanup
# Adds new "scraper" configuration columns (D:J) to the vendor list sheet,
# filling header labels in row 1 and repeating config values for each vendor row (2-5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns D..J (row 1)
$headers = @("searchID", "searchWrapper", "searchItem", "imageWrapper", "imageSrc", "descWrapper", "descSrc")

# New scraper-config values repeated for every data row (2..5), columns D..J
$values = @("main", "link-wrapper", "view-card select-item", "image-wrapper mobile", "data-src", "description", "p")

# Write header row (row 1), columns D(4) .. J(10)
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 4 + $i
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $headers[$i]
    $cell.NumberFormat = "@"
}

# Write the same scraper-config values into every data row (2..5), columns D..J
for ($row = 2; $row -le 5; $row++) {
    for ($i = 0; $i -lt $values.Length; $i++) {
        $col = 4 + $i
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $values[$i]
        $cell.NumberFormat = "@"
    }
}

# Auto-fit the (now widened) columns C through J to their content, like the author did
$ws.Range("C1:J5").EntireColumn.AutoFit()

# Match the author's final selection in the sheet
$ws.Range("D3:J5").Select()
